$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B3").Value = "Mustafa Al Qassimi"
$wsSummary.Range("B4").Value = 1862.57
$wsSummary.Range("B6").Value = 560491
$wsSummary.Range("B7").Value = 372064
$wsSummary.Range("B8").Value = 188427
$wsSummary.Range("B9").Value = 1.51

# ---------------------------------------------------------------------------
# Sheet 2: Assets
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

# Update first data row value
$wsAssets.Range("C2").Value = 145869

# Insert a new row 3 (Vehicles / Luxury Car / 413531), pushing old rows 3-4 down
$wsAssets.Rows.Item(3).Insert()
$wsAssets.Range("A2:C2").Copy()
$wsAssets.Range("A3:C3").PasteSpecial(-4122)
$wsAssets.Range("A3").Value = "Vehicles"
$wsAssets.Range("B3").Value = "Luxury Car"
$wsAssets.Range("C3").Value = 413531

# Update the (now shifted) Liquid Assets row (row 4)
$wsAssets.Range("C4").Value = 1091

# Update the (now shifted) TOTAL ASSETS row (row 5)
$wsAssets.Range("C5").Value = 560491

# ---------------------------------------------------------------------------
# Sheet 3: Liabilities
# ---------------------------------------------------------------------------
$wsLiabilities = $wb.Worksheets.Item("Liabilities")

# Update first data row values
$wsLiabilities.Range("C2").Value = 87521
$wsLiabilities.Range("D2").Value = 1459
$wsLiabilities.Range("E2").Value = 5

# Insert a new row 3 (Auto Loans / Vehicle Loan 2 / 248119 / 6892 / 3), pushing old rows 3-4 down
$wsLiabilities.Rows.Item(3).Insert()
$wsLiabilities.Range("A2:E2").Copy()
$wsLiabilities.Range("A3:E3").PasteSpecial(-4122)
$wsLiabilities.Range("A3").Value = "Auto Loans"
$wsLiabilities.Range("B3").Value = "Vehicle Loan 2"
$wsLiabilities.Range("C3").Value = 248119
$wsLiabilities.Range("D3").Value = 6892
$wsLiabilities.Range("E3").Value = 3

# Update the (now shifted) Credit Cards row (row 4)
$wsLiabilities.Range("C4").Value = 36424
$wsLiabilities.Range("D4").Value = 1821
$wsLiabilities.Range("E4").Value = 1

# Update the (now shifted) TOTAL LIABILITIES row (row 5)
$wsLiabilities.Range("C5").Value = 372064
